# Replace the four "Perseus havainnointijaksot vuonna 2018: ..." observing
# period paragraphs with the new translated Pegasus date text. Each
# paragraph currently contains several runs (and, in one case, a
# hyperlink) that must all collapse into a single plain run with no
# explicit run formatting, matching the target OOXML.

$d = $word.ActiveDocument

$newText = "havainnointijaksot vuonna Pegasus: 8.-17. Lokakuuta 7.-16. Lokakuuta, 7.-16."

$xml = '<?xml version="1.0" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body><w:p><w:r><w:t>' + $newText + '</w:t></w:r></w:p></w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

# Collect the target paragraph indices first (mutating while iterating the
# live collection is unsafe), then apply the replacement to each.
$targets = @()
$i = 0
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Perseus havainnointijaksot*") {
        $targets += $i
    }
    $i++
}

foreach ($idx in $targets) {
    $p = $d.Paragraphs.Item($idx + 1)
    $r = $d.Range($p.Range.Start, $p.Range.End)
    $r.InsertXML($xml)
}

Write-Host "Replaced" $targets.Count "paragraph(s)"
